$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.9711167655286204
$ws.Cells.Item(2, 3).Value = 0.4215847273774216
$ws.Cells.Item(2, 4).Value = 0.02602527397905163
$ws.Cells.Item(2, 6).Value = 0.4650413404363576
$ws.Cells.Item(2, 7).Value = 0.3089105744678378
$ws.Cells.Item(2, 8).Value = 0.4839417890378002
$ws.Cells.Item(2, 9).Value = 0.4865250624665105
$ws.Cells.Item(2, 12).Value = 0.2867844065969294
$ws.Cells.Item(2, 13).Value = 0.2233287334207716
$ws.Cells.Item(2, 15).Value = 1.514733860554671
$ws.Cells.Item(3, 2).Value = 0.8585379497432086
$ws.Cells.Item(3, 3).Value = 0.4071738916452148
$ws.Cells.Item(3, 4).Value = 0.02317560934999108
$ws.Cells.Item(3, 6).Value = 0.4657801849399306
$ws.Cells.Item(3, 7).Value = 0.3108882590470827
$ws.Cells.Item(3, 8).Value = 0.4887947014537701
$ws.Cells.Item(3, 9).Value = 0.4957867113249055
$ws.Cells.Item(3, 12).Value = 0.2860959356829724
$ws.Cells.Item(3, 13).Value = 0.2061702657539684
$ws.Cells.Item(3, 15).Value = 1.528828900888456
$ws.Cells.Item(4, 2).Value = 0.7892128462782466
$ws.Cells.Item(4, 3).Value = 0.3983424049376367
$ws.Cells.Item(4, 4).Value = 0.02141627138114188
$ws.Cells.Item(4, 6).Value = 0.4666282428029263
$ws.Cells.Item(4, 7).Value = 0.3124268792460754
$ws.Cells.Item(4, 8).Value = 0.4920569743757142
$ws.Cells.Item(4, 9).Value = 0.5018498779597405
$ws.Cells.Item(4, 12).Value = 0.2858436360644419
$ws.Cells.Item(4, 13).Value = 0.1956565308294671
$ws.Cells.Item(4, 15).Value = 1.538753444722502
$ws.Cells.Item(5, 2).Value = 0.7609137547648288
$ws.Cells.Item(5, 3).Value = 0.3947481216227402
$ws.Cells.Item(5, 4).Value = 0.0206969508088406
$ws.Cells.Item(5, 6).Value = 0.4670729077802136
$ws.Cells.Item(5, 7).Value = 0.3131352400340361
$ws.Cells.Item(5, 8).Value = 0.4934574235412725
$ws.Cells.Item(5, 9).Value = 0.504415215520039
$ws.Cells.Item(5, 12).Value = 0.2857837164307142
$ws.Cells.Item(5, 13).Value = 0.1913778382494584
$ws.Cells.Item(5, 15).Value = 1.543116737242968
$ws.Cells.Item(6, 2).Value = 0.7562118407688558
$ws.Cells.Item(6, 3).Value = 0.3941515840874672
$ws.Cells.Item(6, 4).Value = 0.02057736602474591
$ws.Cells.Item(6, 6).Value = 0.4671527249163461
$ws.Cells.Item(6, 7).Value = 0.3132577711641389
$ws.Cells.Item(6, 8).Value = 0.4936942578864247
$ws.Cells.Item(6, 9).Value = 0.5048468947947455
$ws.Cells.Item(6, 12).Value = 0.2857763587196303
$ws.Cells.Item(6, 13).Value = 0.1906677201297455
$ws.Cells.Item(6, 15).Value = 1.543860510552392
$ws.Cells.Item(7, 2).Value = 0.7888313886690526
$ws.Cells.Item(7, 3).Value = 0.3982939119562445
$ws.Cells.Item(7, 4).Value = 0.02140657993112427
$ws.Cells.Item(7, 6).Value = 0.46663383869695
$ws.Cells.Item(7, 7).Value = 0.31243610327531
$ws.Cells.Item(7, 8).Value = 0.4920755736677975
$ws.Cells.Item(7, 9).Value = 0.5018840923379475
$ws.Cells.Item(7, 12).Value = 0.2858426542354877
$ws.Cells.Item(7, 13).Value = 0.1955988032462344
$ws.Cells.Item(7, 15).Value = 1.538810998719029
$ws.Cells.Item(8, 2).Value = 0.9323425572241604
$ws.Cells.Item(8, 3).Value = 0.4166126067994185
$ws.Cells.Item(8, 4).Value = 0.02504473559076104
$ws.Cells.Item(8, 6).Value = 0.465214173155502
$ws.Cells.Item(8, 7).Value = 0.3095250537716794
$ws.Cells.Item(8, 8).Value = 0.4855564228452565
$ws.Cells.Item(8, 9).Value = 0.4896402701096694
$ws.Cells.Item(8, 12).Value = 0.2865116731663946
$ws.Cells.Item(8, 13).Value = 0.217408177686508
$ws.Cells.Item(8, 15).Value = 1.519329958779238
$ws.Cells.Item(9, 2).Value = 1.212093938135752
$ws.Cells.Item(9, 3).Value = 0.4526538308719523
$ws.Cells.Item(9, 4).Value = 0.03210101992269898
$ws.Cells.Item(9, 6).Value = 0.4655649192871394
$ws.Cells.Item(9, 7).Value = 0.3063986875728233
$ws.Cells.Item(9, 8).Value = 0.4750147928833002
$ws.Cells.Item(9, 9).Value = 0.4686224144566857
$ws.Cells.Item(9, 12).Value = 0.2891745862570261
$ws.Cells.Item(9, 13).Value = 0.260336970606545
$ws.Cells.Item(9, 15).Value = 1.491225573793471
$ws.Cells.Item(10, 2).Value = 1.416524092291638
$ws.Cells.Item(10, 3).Value = 0.4791878283783149
$ws.Cells.Item(10, 4).Value = 0.03723583726576862
$ws.Cells.Item(10, 6).Value = 0.4677419344235076
$ws.Cells.Item(10, 7).Value = 0.3056886828947967
$ws.Cells.Item(10, 8).Value = 0.4686375143487993
$ws.Cells.Item(10, 9).Value = 0.4550112516556748
$ws.Cells.Item(10, 12).Value = 0.2919536218731622
$ws.Cells.Item(10, 13).Value = 0.2919631864818939
$ws.Cells.Item(10, 15).Value = 1.476761965364162
$ws.Cells.Item(11, 2).Value = 1.50926894533859
$ws.Cells.Item(11, 3).Value = 0.4912671437463132
$ws.Cells.Item(11, 4).Value = 0.03956071096725111
$ws.Cells.Item(11, 6).Value = 0.469150902069984
$ws.Cells.Item(11, 7).Value = 0.3057130146031142
$ws.Cells.Item(11, 8).Value = 0.4660334442884562
$ws.Cells.Item(11, 9).Value = 0.4492180602520204
$ws.Cells.Item(11, 12).Value = 0.2933963296196964
$ws.Cells.Item(11, 13).Value = 0.3063673138310321
$ws.Cells.Item(11, 15).Value = 1.471531240661051
$ws.Cells.Item(12, 2).Value = 1.544351068179367
$ws.Cells.Item(12, 3).Value = 0.4958421924705192
$ws.Cells.Item(12, 4).Value = 0.04043946063588066
$ws.Cells.Item(12, 6).Value = 0.4697447643540542
$ws.Cells.Item(12, 7).Value = 0.3057723819986151
$ws.Cells.Item(12, 8).Value = 0.4650900702417147
$ws.Cells.Item(12, 9).Value = 0.447081773282985
$ws.Cells.Item(12, 12).Value = 0.2939682866448408
$ws.Cells.Item(12, 13).Value = 0.3118239997381096
$ws.Cells.Item(12, 15).Value = 1.469744942359512
$ws.Cells.Item(13, 2).Value = 1.536797245104083
$ws.Cells.Item(13, 3).Value = 0.494856841339157
$ws.Cells.Item(13, 4).Value = 0.04025027944967974
$ws.Cells.Item(13, 6).Value = 0.4696141810256265
$ws.Cells.Item(13, 7).Value = 0.3057573622594489
$ws.Cells.Item(13, 8).Value = 0.4652913418733533
$ws.Cells.Item(13, 9).Value = 0.4475393025572263
$ws.Cells.Item(13, 12).Value = 0.2938439659023828
$ws.Cells.Item(13, 13).Value = 0.3106487139801999
$ws.Cells.Item(13, 15).Value = 1.47012099691085
$ws.Cells.Item(14, 2).Value = 1.512155953540628
$ws.Cells.Item(14, 3).Value = 0.4916435205225582
$ws.Cells.Item(14, 4).Value = 0.03963303917851135
$ws.Cells.Item(14, 6).Value = 0.4691985499712885
$ws.Cells.Item(14, 7).Value = 0.305716892449432
$ws.Cells.Item(14, 8).Value = 0.4659549758105044
$ws.Cells.Item(14, 9).Value = 0.4490411540103985
$ws.Cells.Item(14, 12).Value = 0.2934428713063824
$ws.Cells.Item(14, 13).Value = 0.3068161972828776
$ws.Cells.Item(14, 15).Value = 1.471380380630876
$ws.Cells.Item(15, 2).Value = 1.497057393419766
$ws.Cells.Item(15, 3).Value = 0.4896753707392065
$ws.Cells.Item(15, 4).Value = 0.03925474797657102
$ws.Cells.Item(15, 6).Value = 0.4689518224404736
$ws.Cells.Item(15, 7).Value = 0.3056986411790774
$ws.Cells.Item(15, 8).Value = 0.46636703628851
$ws.Cells.Item(15, 9).Value = 0.4499685700179992
$ws.Cells.Item(15, 12).Value = 0.2932005266729192
$ws.Cells.Item(15, 13).Value = 0.3044689414336545
$ws.Cells.Item(15, 15).Value = 1.472177129218522
$ws.Cells.Item(16, 2).Value = 1.410457751359843
$ws.Cells.Item(16, 3).Value = 0.4783985607893158
$ws.Cells.Item(16, 4).Value = 0.03708367602193618
$ws.Cells.Item(16, 6).Value = 0.4676582891684333
$ws.Cells.Item(16, 7).Value = 0.3056940996415136
$ws.Cells.Item(16, 8).Value = 0.4688136746756726
$ws.Cells.Item(16, 9).Value = 0.4553978840471622
$ws.Cells.Item(16, 12).Value = 0.2918629266596326
$ws.Cells.Item(16, 13).Value = 0.2910221618087263
$ws.Cells.Item(16, 15).Value = 1.477130990706556
$ws.Cells.Item(17, 2).Value = 1.357265716447785
$ws.Cells.Item(17, 3).Value = 0.4714825878475608
$ws.Cells.Item(17, 4).Value = 0.03574894532148676
$ws.Cells.Item(17, 6).Value = 0.4669720549132137
$ws.Cells.Item(17, 7).Value = 0.3057804196138036
$ws.Cells.Item(17, 8).Value = 0.4703906909054467
$ws.Cells.Item(17, 9).Value = 0.4588307900272728
$ws.Cells.Item(17, 12).Value = 0.2910880486234788
$ws.Cells.Item(17, 13).Value = 0.2827771771876897
$ws.Cells.Item(17, 15).Value = 1.480515854073943
$ws.Cells.Item(18, 2).Value = 1.32664752672548
$ws.Cells.Item(18, 3).Value = 0.4675055614675045
$ws.Cells.Item(18, 4).Value = 0.03498021300656262
$ws.Cells.Item(18, 6).Value = 0.4666167478914005
$ws.Cells.Item(18, 7).Value = 0.3058627472375832
$ws.Cells.Item(18, 8).Value = 0.4713257012092598
$ws.Cells.Item(18, 9).Value = 0.4608428153740345
$ws.Cells.Item(18, 12).Value = 0.2906591636210578
$ws.Cells.Item(18, 13).Value = 0.2780365182665676
$ws.Cells.Item(18, 15).Value = 1.482589667839349
$ws.Cells.Item(19, 2).Value = 1.316276757814535
$ws.Cells.Item(19, 3).Value = 0.4661591683619406
$ws.Cells.Item(19, 4).Value = 0.0347197580037033
$ws.Cells.Item(19, 6).Value = 0.4665032098403614
$ws.Cells.Item(19, 7).Value = 0.3058962281187689
$ws.Cells.Item(19, 8).Value = 0.4716470797095624
$ws.Cells.Item(19, 9).Value = 0.4615304901019108
$ws.Cells.Item(19, 12).Value = 0.2905168377426293
$ws.Cells.Item(19, 13).Value = 0.2764317036929285
$ws.Cells.Item(19, 15).Value = 1.483313610504325
$ws.Cells.Item(20, 2).Value = 1.362930551329327
$ws.Cells.Item(20, 3).Value = 0.4722187190520799
$ws.Cells.Item(20, 4).Value = 0.03589113670463462
$ws.Cells.Item(20, 6).Value = 0.4670410274413825
$ws.Cells.Item(20, 7).Value = 0.3057678470602667
$ws.Cells.Item(20, 8).Value = 0.4702199217003198
$ws.Cells.Item(20, 9).Value = 0.4584614682239145
$ws.Cells.Item(20, 12).Value = 0.2911687967770575
$ws.Cells.Item(20, 13).Value = 0.283654702643382
$ws.Cells.Item(20, 15).Value = 1.480142389171732
$ws.Cells.Item(21, 2).Value = 1.519394756045301
$ws.Cells.Item(21, 3).Value = 0.4925873301916113
$ws.Cells.Item(21, 4).Value = 0.03981438223617317
$ws.Cells.Item(21, 6).Value = 0.4693189929774348
$ws.Cells.Item(21, 7).Value = 0.3057274166005755
$ws.Cells.Item(21, 8).Value = 0.4657588905419416
$ws.Cells.Item(21, 9).Value = 0.4485984629409732
$ws.Cells.Item(21, 12).Value = 0.2935599870659189
$ws.Cells.Item(21, 13).Value = 0.3079418443676403
$ws.Cells.Item(21, 15).Value = 1.471005187447872
$ws.Cells.Item(22, 2).Value = 1.621428680451231
$ws.Cells.Item(22, 3).Value = 0.5059043065663218
$ws.Cells.Item(22, 4).Value = 0.04236892836860306
$ws.Cells.Item(22, 6).Value = 0.4711594007894107
$ws.Cells.Item(22, 7).Value = 0.3059934217837679
$ws.Cells.Item(22, 8).Value = 0.4630924275202233
$ws.Cells.Item(22, 9).Value = 0.4424874612360288
$ws.Cells.Item(22, 12).Value = 0.2952721717146432
$ws.Cells.Item(22, 13).Value = 0.3238273127555757
$ws.Cells.Item(22, 15).Value = 1.466167209536593
$ws.Cells.Item(23, 2).Value = 1.566992511701244
$ws.Cells.Item(23, 3).Value = 0.4987964678449828
$ws.Cells.Item(23, 4).Value = 0.04100640692719537
$ws.Cells.Item(23, 6).Value = 0.4701449283597725
$ws.Cells.Item(23, 7).Value = 0.3058246237949902
$ws.Cells.Item(23, 8).Value = 0.464492770294072
$ws.Cells.Item(23, 9).Value = 0.445718311606603
$ws.Cells.Item(23, 12).Value = 0.2943446876494846
$ws.Cells.Item(23, 13).Value = 0.3153479078801169
$ws.Cells.Item(23, 15).Value = 1.468645430427244
$ws.Cells.Item(24, 2).Value = 1.360369595749773
$ws.Cells.Item(24, 3).Value = 0.4718859170368717
$ws.Cells.Item(24, 4).Value = 0.03582685626818716
$ws.Cells.Item(24, 6).Value = 0.4670097228093368
$ws.Cells.Item(24, 7).Value = 0.3057734292668428
$ws.Cells.Item(24, 8).Value = 0.4702970381120721
$ws.Cells.Item(24, 9).Value = 0.4586283190033349
$ws.Cells.Item(24, 12).Value = 0.2911322388194577
$ws.Cells.Item(24, 13).Value = 0.2832579749416624
$ws.Cells.Item(24, 15).Value = 1.480310834571483
$ws.Cells.Item(25, 2).Value = 1.136601576579039
$ws.Cells.Item(25, 3).Value = 0.4428928923975093
$ws.Cells.Item(25, 4).Value = 0.03020065935703542
$ws.Cells.Item(25, 6).Value = 0.4651335732093855
$ws.Cells.Item(25, 7).Value = 0.3069666841586738
$ws.Cells.Item(25, 8).Value = 0.4776264521702842
$ws.Cells.Item(25, 9).Value = 0.4739872866007673
$ws.Cells.Item(25, 12).Value = 0.2883096689441089
$ws.Cells.Item(25, 13).Value = 0.2487075908363323
$ws.Cells.Item(25, 15).Value = 1.497744517184714
